$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column AD ("area") holds value 10 for rows 2 through 31.
# Update it to 8 for all those rows (smaller, enlarged areas per commit message).
$ws.Range("AD2:AD31").Value = 8
